$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-1651255574212943"
$ws1.Range("B2").Value = "go_stims-16512555741813917.csv"
$ws1.Range("B3").Value = "GNG_stims-1651255574198038.csv"
$ws1.Range("B4").Value = "go_stims-16512555741990266.csv"
$ws1.Range("B5").Value = "GNG_stims-1651255574211944.csv"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16512555762625153"
$ws2.Range("B2").Value = "OB-1651255575524035.csv"
$ws2.Range("B3").Value = "ZB-match_1-16512555743809593.csv"
$ws2.Range("B4").Value = "ZB-match_4-1651255574794492.csv"
$ws2.Range("B5").Value = "TB-16512555758231156.csv"
$ws2.Range("B6").Value = "TB-1651255575773884.csv"
$ws2.Range("B7").Value = "OB-1651255575738886.csv"
$ws2.Range("B8").Value = "ZB-match_9-16512555746536973.csv"
$ws2.Range("B9").Value = "TB-1651255576231916.csv"
$ws2.Range("B10").Value = "OB-16512555753776557.csv"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1651255576264542"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16512555763100703"
$ws4.Range("B2").Value = "MM_stims-1651255576277993.csv"
$ws4.Range("B3").Value = "ZM_stims-1651255576266543.csv"
$ws4.Range("B4").Value = "MM_stims-16512555762940729.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555762789967.csv"
$ws4.Range("B6").Value = "MM_stims-16512555763100703.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555762951312.csv"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16512555763894393"
$ws5.Range("B2").Value = "vSAT_stims-16512555763584886.csv"
$ws5.Range("B3").Value = "vSAT_stims-16512555763734372.csv"
$ws5.Range("B4").Value = "SAT_stims-16512555763434258.csv"
$ws5.Range("B5").Value = "SAT_stims-16512555763170507.csv"
